$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the seven Disney-licensed rows (IW_093 .. IW_099) that are being
# dropped from the data model.
$ws.Rows("94:100").Delete()

# The remaining rows (IW_001 .. IW_092) get a new "License List" value in
# column D and the "Image Directory" value in column E is repointed from the
# old Image_Wonderland folder to the new shared Image folder.
# NOTE: write column E's new shared string before column D's so the shared
# string table ends up ordered the same way the workbook author's edit did.
for ($r = 2; $r -le 93; $r++) {
    $ws.Cells.Item($r, 5).Value = "data/Multimedia_Data/Image/"
    $ws.Cells.Item($r, 4).Value = "CC 0.0"
}

# Re-point the view: scroll so column A is visible again and select the new
# License List column that was just populated.
$ws.Range("D2:D93").Select()
